$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7573912063605235
$ws.Range("C2").Value = 0.1350596674196183
$ws.Range("D2").Value = 0.07766267278526584
$ws.Range("E2").Value = 0.1114702786803434
$ws.Range("G2").Value = 0.002503515525458773
$ws.Range("K2").Value = 0.3964512192575569
$ws.Range("L2").Value = 0.1923488513829525
$ws.Range("M2").Value = 0.1964889248670012
$ws.Range("N2").Value = 2.414080146706652
$ws.Range("O2").Value = 5.077439800590327
$ws.Range("B3").Value = 0.7215847792689374
$ws.Range("C3").Value = 0.1334551369442636
$ws.Range("D3").Value = 0.07057659360374657
$ws.Range("E3").Value = 0.1117850815743875
$ws.Range("G3").Value = 0.002506345585331836
$ws.Range("K3").Value = 0.3624179940076431
$ws.Range("L3").Value = 0.1898298979802036
$ws.Range("M3").Value = 0.189943701298553
$ws.Range("N3").Value = 2.434714652742024
$ws.Range("O3").Value = 5.087404786434831
$ws.Range("B4").Value = 0.6999482295969699
$ws.Range("C4").Value = 0.1324539530591196
$ws.Range("D4").Value = 0.06626029019173529
$ws.Range("E4").Value = 0.1120163324784844
$ws.Range("G4").Value = 0.002508177253422235
$ws.Range("K4").Value = 0.3416498643481702
$ws.Range("L4").Value = 0.1883750849257737
$ws.Range("M4").Value = 0.1860203510491552
$ws.Range("N4").Value = 2.448069087572897
$ws.Range("O4").Value = 5.095807171551911
$ws.Range("B5").Value = 0.6912193450832262
$ws.Range("C5").Value = 0.1320419454194095
$ws.Range("D5").Value = 0.06451005918691521
$ws.Range("E5").Value = 0.1121201318774201
$ws.Range("G5").Value = 0.002508947380574696
$ws.Range("K5").Value = 0.3332193132867474
$ws.Range("L5").Value = 0.1878053835403364
$ws.Range("M5").Value = 0.1844456522146345
$ws.Range("N5").Value = 2.453683337322747
$ws.Range("O5").Value = 5.099805781788177
$ws.Range("B6").Value = 0.6897752611117767
$ws.Range("C6").Value = 0.1319732893969388
$ws.Range("D6").Value = 0.06421996006685049
$ws.Range("E6").Value = 0.1121379457651024
$ws.Range("G6").Value = 0.002509076693655526
$ws.Range("K6").Value = 0.3318214075657266
$ws.Range("L6").Value = 0.1877121846638872
$ws.Range("M6").Value = 0.1841856330430183
$ws.Range("N6").Value = 2.454625980840262
$ws.Range("O6").Value = 5.10050445889658
$ws.Range("B7").Value = 0.6998301510976148
$ws.Range("C7").Value = 0.1324484128386061
$ws.Range("D7").Value = 0.06623665074405949
$ws.Range("E7").Value = 0.1120176936110848
$ws.Range("G7").Value = 0.002508187543580809
$ws.Range("K7").Value = 0.3415360344164355
$ws.Range("L7").Value = 0.1883673079412418
$ws.Range("M7").Value = 0.1859990164066225
$ws.Range("N7").Value = 2.448144106017477
$ws.Range("O7").Value = 5.095858771432205
$ws.Range("B8").Value = 0.7449730509968902
$ws.Range("C8").Value = 0.1345097479094832
$ws.Range("D8").Value = 0.07521221556714863
$ws.Range("E8").Value = 0.1115709542200083
$ws.Range("G8").Value = 0.002504471862995376
$ws.Range("K8").Value = 0.3846901057019068
$ws.Range("L8").Value = 0.1914612872642323
$ws.Range("M8").Value = 0.1942123789022396
$ws.Range("N8").Value = 2.421052814905053
$ws.Range("O8").Value = 5.080401873549533
$ws.Range("B9").Value = 0.836248852546305
$ws.Range("C9").Value = 0.1384250343431148
$ws.Range("D9").Value = 0.09308845118829367
$ws.Range("E9").Value = 0.1109953933686754
$ws.Range("G9").Value = 0.002497928043367631
$ws.Range("K9").Value = 0.4703234439374171
$ws.Range("L9").Value = 0.1982553803847722
$ws.Range("M9").Value = 0.2110727682288527
$ws.Range("N9").Value = 2.373358585190935
$ws.Range("O9").Value = 5.068204426831329
$ws.Range("B10").Value = 0.9049717218957483
$ws.Range("C10").Value = 0.1412243287063646
$ws.Range("D10").Value = 0.1063923782993328
$ws.Range("E10").Value = 0.1107548432701773
$ws.Range("G10").Value = 0.002493568498157012
$ws.Range("K10").Value = 0.5338454047399068
$ws.Range("L10").Value = 0.2036883630949404
$ws.Range("M10").Value = 0.2239168501303297
$ws.Range("N10").Value = 2.341625729545889
$ws.Range("O10").Value = 5.070281586121155
$ws.Range("B11").Value = 0.9365939610424334
$ws.Range("C11").Value = 0.1424810677905484
$ws.Range("D11").Value = 0.1124822568681338
$ws.Range("E11").Value = 0.1106848194277958
$ws.Range("G11").Value = 0.002491681589262099
$ws.Range("K11").Value = 0.5628738184286703
$ws.Range("L11").Value = 0.2062554898103599
$ws.Range("M11").Value = 0.2298586112214096
$ws.Range("N11").Value = 2.327907279343151
$ws.Range("O11").Value = 5.073623190565172
$ws.Range("B12").Value = 0.9486198382041948
$ws.Range("C12").Value = 0.1429545624220552
$ws.Range("D12").Value = 0.114793795482683
$ws.Range("E12").Value = 0.1106639535110094
$ws.Range("G12").Value = 0.002490980836719899
$ws.Range("K12").Value = 0.5738848437752893
$ws.Range("L12").Value = 0.2072413031869758
$ws.Range("M12").Value = 0.2321227489364972
$ws.Range("N12").Value = 2.322815552042833
$ws.Range("O12").Value = 5.075233080810392
$ws.Range("B13").Value = 0.9460275806774519
$ws.Range("C13").Value = 0.1428526939715908
$ws.Range("D13").Value = 0.1142957227780244
$ws.Range("E13").Value = 0.1106681962958351
$ws.Range("G13").Value = 0.002491131144698487
$ws.Range("K13").Value = 0.5715126014871998
$ws.Range("L13").Value = 0.20702838226849
$ws.Range("M13").Value = 0.231634499947539
$ws.Range("N13").Value = 2.3239075580668
$ws.Range("O13").Value = 5.074871042814038
$ws.Range("B14").Value = 0.9375823136526265
$ws.Range("C14").Value = 0.1425200707250482
$ws.Range("D14").Value = 0.1126723196983477
$ws.Range("E14").Value = 0.1106829896190984
$ws.Range("G14").Value = 0.002491623662258421
$ws.Range("K14").Value = 0.5637793316152795
$ws.Range("L14").Value = 0.2063363190476935
$ws.Range("M14").Value = 0.2300446007655381
$ws.Range("N14").Value = 2.327486312253036
$ws.Range("O14").Value = 5.073748733788818
$ws.Range("B15").Value = 0.9324159986799714
$ws.Range("C15").Value = 0.1423160159341208
$ws.Range("D15").Value = 0.1116786457720877
$ws.Range("E15").Value = 0.1106927863686913
$ws.Range("G15").Value = 0.002491927136015643
$ws.Range("K15").Value = 0.5590448908706662
$ws.Range("L15").Value = 0.2059141926463468
$ws.Range("M15").Value = 0.2290725776074893
$ws.Range("N15").Value = 2.32969183752008
$ws.Range("O15").Value = 5.073106146552078
$ws.Range("B16").Value = 0.9029123305586779
$ws.Range("C16").Value = 0.1411418616883822
$ws.Range("D16").Value = 0.1059951513391155
$ws.Range("E16").Value = 0.1107602110597821
$ws.Range("G16").Value = 0.002493693743930583
$ws.Range("K16").Value = 0.5319509524467492
$ws.Range("L16").Value = 0.203522515137692
$ws.Range("M16").Value = 0.2235305238336878
$ws.Range("N16").Value = 2.342536694346862
$ws.Range("O16").Value = 5.070111419619252
$ws.Range("B17").Value = 0.8849045806162508
$ws.Range("C17").Value = 0.140417280156349
$ws.Range("D17").Value = 0.1025182005220984
$ws.Range("E17").Value = 0.1108116550184466
$ws.Range("G17").Value = 0.002494802111434818
$ws.Range("K17").Value = 0.5153632255118055
$ws.Range("L17").Value = 0.2020797545631865
$ws.Range("M17").Value = 0.2201559165772622
$ws.Range("N17").Value = 2.350600269986018
$ws.Range("O17").Value = 5.068887998993262
$ws.Range("B18").Value = 0.8745809024402149
$ws.Range("C18").Value = 0.1399989510403259
$ws.Range("D18").Value = 0.1005219141351148
$ws.Range("E18").Value = 0.1108449540065202
$ws.Range("G18").Value = 0.002495448680500527
$ws.Range("K18").Value = 0.5058348466884866
$ws.Range("L18").Value = 0.201258921755894
$ws.Range("M18").Value = 0.2182242491007855
$ws.Range("N18").Value = 2.355305708211418
$ws.Range("O18").Value = 5.068409916480562
$ws.Range("B19").Value = 0.8710913213476488
$ws.Range("C19").Value = 0.1398570427427046
$ws.Range("D19").Value = 0.099846618230643
$ws.Range("E19").Value = 0.110856866106511
$ws.Range("G19").Value = 0.002495669156436161
$ws.Range("K19").Value = 0.5026108500007922
$ws.Range("L19").Value = 0.2009825501210258
$ws.Range("M19").Value = 0.2175718230987229
$ws.Range("N19").Value = 2.356910477496157
$ws.Range("O19").Value = 5.068286796857905
$ws.Range("B20").Value = 0.8868180312846334
$ws.Range("C20").Value = 0.1404945755356408
$ws.Range("D20").Value = 0.1028879591057859
$ws.Range("E20").Value = 0.1108057948858168
$ws.Range("G20").Value = 0.002494683185953727
$ws.Range("K20").Value = 0.5171277323107688
$ws.Range("L20").Value = 0.2022324073067239
$ws.Range("M20").Value = 0.2205141857174411
$ws.Range("N20").Value = 2.349734904016628
$ws.Range("O20").Value = 5.068994886843115
$ws.Range("B21").Value = 0.9400615073523682
$ws.Range("C21").Value = 0.1426178355707748
$ws.Range("D21").Value = 0.1131490051457007
$ws.Range("E21").Value = 0.1106784912301588
$ws.Range("G21").Value = 0.002491478624358561
$ws.Range("K21").Value = 0.5660502787989401
$ws.Range("L21").Value = 0.2065392234575256
$ws.Range("M21").Value = 0.2305112102592801
$ws.Range("N21").Value = 2.32643234546326
$ws.Range("O21").Value = 5.074069034741257
$ws.Range("B22").Value = 0.9751575105671577
$ws.Range("C22").Value = 0.1439914860110107
$ws.Range("D22").Value = 0.1198868590244189
$ws.Range("E22").Value = 0.110628219705692
$ws.Range("G22").Value = 0.002489464542836939
$ws.Range("K22").Value = 0.5981322007636436
$ws.Range("L22").Value = 0.2094337952642462
$ws.Range("M22").Value = 0.2371271080035839
$ws.Range("N22").Value = 2.311804051123552
$ws.Range("O22").Value = 5.079393211347082
$ws.Range("B23").Value = 0.9563989913210946
$ws.Range("C23").Value = 0.1432596284736576
$ws.Range("D23").Value = 0.1162878454744316
$ws.Range("E23").Value = 0.1106520427202788
$ws.Range("G23").Value = 0.002490532171169155
$ws.Range("K23").Value = 0.5809996996828772
$ws.Range("L23").Value = 0.2078816222989985
$ws.Range("M23").Value = 0.2335885858541005
$ws.Range("N23").Value = 2.31955641297138
$ws.Range("O23").Value = 5.076367919773105
$ws.Range("B24").Value = 0.8859528691903051
$ws.Range("C24").Value = 0.1404596357693464
$ws.Range("D24").Value = 0.1027207829683476
$ws.Range("E24").Value = 0.1108084326535721
$ws.Range("G24").Value = 0.002494736923195682
$ws.Range("K24").Value = 0.5163299734468012
$ws.Range("L24").Value = 0.2021633661180005
$ws.Range("M24").Value = 0.2203521859445488
$ws.Range("N24").Value = 2.350125919225981
$ws.Range("O24").Value = 5.06894586109496
$ws.Range("B25").Value = 0.8112633326279877
$ws.Range("C25").Value = 0.1373794250997378
$ws.Range("D25").Value = 0.08822273510917
$ws.Range("E25").Value = 0.1111190185221513
$ws.Range("G25").Value = 0.002499619286660152
$ws.Range("K25").Value = 0.4470501883253348
$ws.Range("L25").Value = 0.1963397660244866
$ws.Range("M25").Value = 0.2110727682288527
$ws.Range("N25").Value = 2.385679981190123
$ws.Range("O25").Value = 5.069565394598982
